$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.974.81"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.138.35"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.90"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.68"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.135.76"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("E10").Value = "  +5.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.503"
$ws.Range("E12").Value = "  +6.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  +11.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.36"
$ws.Range("E14").Value = "  +6.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.652.28"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.982.74"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("E17").Value = "  +6.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.143.06"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.24"
$ws.Range("E20").Value = "  +6.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.92"
$ws.Range("E21").Value = "  +7.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").Value = "  +8.45%  "
$ws.Range("E23").Value = "  +12.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.83"
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.20"
$ws.Range("E25").Value = "  +4.98%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("E28").Value = "  +8.43%  "
$ws.Range("E29").Value = "  +5.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.87"
$ws.Range("E30").Value = "  +6.47%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.65"
$ws.Range("E33").Value = "  +6.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.02"
$ws.Range("E34").Value = "  +8.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.58"
$ws.Range("E35").Value = "  +6.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.65"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "475.46"
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0858"
$ws.Range("E39").Value = "  +4.06%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.106.77"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.61"
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.120"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("E44").Value = "  +10.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  +13.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.14"
$ws.Range("E46").Value = "  +4.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0577"
$ws.Range("E47").Value = "  +11.74%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +3.53%  "
$ws.Range("E50").Value = "  +10.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.95"
$ws.Range("E51").Value = "  -1.77%  "
